$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide rows that were previously hidden (several new factors became candidates) ---
$rowsToUnhide = @(3,4,5,9,10,11,12,15,16,17,18,20,22,23)
foreach ($r in $rowsToUnhide) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- Update "Do I need?" (col C) / "Done?" (col D) values for the affected factors ---
# Row 6 - Credit Rating: now sourced from TRACE, possible to get
$ws.Range("D6").Value = "Possible"
$ws.Range("C6").Value = "TRACE"

# Row 19 - Duration*: not needed, marked done
$ws.Range("D19").Value = "Done"
$ws.Range("C19").Value = "No"

# Row 22 - Default-beta
$ws.Range("D22").Value = "Possible"

# Row 23 - Term-beta
$ws.Range("D23").Value = "Possible"

# Row 24 - Yield to maturity
$ws.Range("D24").Value = "Done"

# Row 25 - Credit Spread
$ws.Range("D25").Value = "Possible"

# Row 26 - Bond Age
$ws.Range("D26").Value = "Done"
$ws.Range("C26").Value = "TRACE"

# Row 27 - Amount Outstanding*
$ws.Range("D27").Value = "Done"
$ws.Range("C27").Value = "TRACE"

# Row 28 - Market Value*
$ws.Range("D28").Value = "Done"
$ws.Range("C28").Value = "No"

# Row 29 - Equity Momentum
$ws.Range("D29").Value = "Possible"

# --- Clear the "Do I need?" column filter on Table1 (kept the filter column, dropped the criteria) ---
$lo = $ws.ListObjects.Item(1)
$tableRange = $lo.Range
$tableName = $lo.Name
$lo.Unlist()
$newlo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$newlo.Name = $tableName
$newlo.TableStyle = "TableStyleMedium1"

# --- Move the active selection to D6 ---
$ws.Range("D6").Select()
